$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.20545
$ws.Range("H2").Value = 12.61635
$ws.Range("I2").Value = 0.0197921807762369
$ws.Range("J2").Value = 0.0197921807762369
$ws.Range("M2").Value = 0.6068583333333333
$ws.Range("N2").Value = 1.820575
$ws.Range("O2").Value = 0.2510732750714712
$ws.Range("P2").Value = 0.2510732750714712
$ws.Range("Q2").Value = 2.552112377916667
$ws.Range("R2").Value = 22.96901140125
$ws.Range("S2").Value = 0.004969287648296411
$ws.Range("T2").Value = 0.004969287648296411
$ws.Range("G3").Value = 4.20545
$ws.Range("H3").Value = 12.61635
$ws.Range("I3").Value = 0.0197921807762369
$ws.Range("J3").Value = 0.0197921807762369
$ws.Range("O3").Value = 0.08210468103768082
$ws.Range("P3").Value = 0.08210468103768082
$ws.Range("Q3").Value = 0.8345785615833333
$ws.Range("R3").Value = 7.51120705425
$ws.Range("S3").Value = 0.001625030689673048
$ws.Range("T3").Value = 0.001625030689673048
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 4.20545
$ws.Range("H4").Value = 12.61635
$ws.Range("I4").Value = 0.0197921807762369
$ws.Range("J4").Value = 0.0197921807762369
$ws.Range("M4").Value = 0.6741056666666667
$ws.Range("N4").Value = 2.022317
$ws.Range("O4").Value = 0.278895267936071
$ws.Range("P4").Value = 0.278895267936071
$ws.Range("Q4").Value = 2.834917675883334
$ws.Range("R4").Value = 25.51425908295
$ws.Range("S4").Value = 0.005519945560627743
$ws.Range("T4").Value = 0.005519945560627743
$ws.Range("G5").Value = 4.20545
$ws.Range("H5").Value = 12.61635
$ws.Range("I5").Value = 0.0197921807762369
$ws.Range("J5").Value = 0.0197921807762369
$ws.Range("M5").Value = 0.9376410000000001
$ws.Range("N5").Value = 2.812923
$ws.Range("O5").Value = 0.3879267759547769
$ws.Range("P5").Value = 0.387926775954777
$ws.Range("Q5").Value = 3.94320234345
$ws.Range("R5").Value = 35.48882109105
$ws.Range("S5").Value = 0.007677916877639693
$ws.Range("T5").Value = 0.007677916877639695
$ws.Range("I6").Value = 0.02862974203518836
$ws.Range("J6").Value = 0.02862974203518836
$ws.Range("M6").Value = 0.6068583333333333
$ws.Range("N6").Value = 1.820575
$ws.Range("O6").Value = 0.2510732750714712
$ws.Range("P6").Value = 0.2510732750714712
$ws.Range("Q6").Value = 3.691676013402777
$ws.Range("R6").Value = 33.225084120625
$ws.Range("S6").Value = 0.007188163097226109
$ws.Range("T6").Value = 0.007188163097226109
$ws.Range("I7").Value = 0.02862974203518836
$ws.Range("J7").Value = 0.02862974203518836
$ws.Range("O7").Value = 0.08210468103768082
$ws.Range("P7").Value = 0.08210468103768082
$ws.Range("S7").Value = 0.002350635837990223
$ws.Range("T7").Value = 0.002350635837990223
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.02862974203518836
$ws.Range("J8").Value = 0.02862974203518836
$ws.Range("M8").Value = 0.6741056666666667
$ws.Range("N8").Value = 2.022317
$ws.Range("O8").Value = 0.278895267936071
$ws.Range("P8").Value = 0.278895267936071
$ws.Range("Q8").Value = 4.100758914297223
$ws.Range("R8").Value = 36.90683022867501
$ws.Range("S8").Value = 0.007984699575844453
$ws.Range("T8").Value = 0.007984699575844453
$ws.Range("I9").Value = 0.02862974203518836
$ws.Range("J9").Value = 0.02862974203518836
$ws.Range("M9").Value = 0.9376410000000001
$ws.Range("N9").Value = 2.812923
$ws.Range("O9").Value = 0.3879267759547769
$ws.Range("P9").Value = 0.387926775954777
$ws.Range("Q9").Value = 5.703912426925
$ws.Range("R9").Value = 51.335211842325
$ws.Range("S9").Value = 0.01110624352412757
$ws.Range("T9").Value = 0.01110624352412758
$ws.Range("G10").Value = 99.54897833333332
$ws.Range("H10").Value = 298.646935
$ws.Range("I10").Value = 0.4685090478457771
$ws.Range("J10").Value = 0.4685090478457771
$ws.Range("M10").Value = 0.6068583333333333
$ws.Range("N10").Value = 1.820575
$ws.Range("O10").Value = 0.2510732750714712
$ws.Range("P10").Value = 0.2510732750714712
$ws.Range("Q10").Value = 60.41212707640277
$ws.Range("R10").Value = 543.7091436876249
$ws.Range("S10").Value = 0.1176301010432558
$ws.Range("T10").Value = 0.1176301010432558
$ws.Range("G11").Value = 99.54897833333332
$ws.Range("H11").Value = 298.646935
$ws.Range("I11").Value = 0.4685090478457771
$ws.Range("J11").Value = 0.4685090478457771
$ws.Range("O11").Value = 0.08210468103768082
$ws.Range("P11").Value = 0.08210468103768082
$ws.Range("Q11").Value = 19.75566066521389
$ws.Range("R11").Value = 177.800945986925
$ws.Range("S11").Value = 0.03846678593664506
$ws.Range("T11").Value = 0.03846678593664506
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 99.54897833333332
$ws.Range("H12").Value = 298.646935
$ws.Range("I12").Value = 0.4685090478457771
$ws.Range("J12").Value = 0.4685090478457771
$ws.Range("M12").Value = 0.6741056666666667
$ws.Range("N12").Value = 2.022317
$ws.Range("O12").Value = 0.278895267936071
$ws.Range("P12").Value = 0.278895267936071
$ws.Range("Q12").Value = 67.10653040537721
$ws.Range("R12").Value = 603.958773648395
$ws.Range("S12").Value = 0.1306649564294215
$ws.Range("T12").Value = 0.1306649564294215
$ws.Range("G13").Value = 99.54897833333332
$ws.Range("H13").Value = 298.646935
$ws.Range("I13").Value = 0.4685090478457771
$ws.Range("J13").Value = 0.4685090478457771
$ws.Range("M13").Value = 0.9376410000000001
$ws.Range("N13").Value = 2.812923
$ws.Range("O13").Value = 0.3879267759547769
$ws.Range("P13").Value = 0.387926775954777
$ws.Range("Q13").Value = 93.341203593445
$ws.Range("R13").Value = 840.070832341005
$ws.Range("S13").Value = 0.1817472044364546
$ws.Range("T13").Value = 0.1817472044364546
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1674636666666667
$ws.Range("H14").Value = 0.502391
$ws.Range("I14").Value = 0.0007881370992683645
$ws.Range("J14").Value = 0.0007881370992683645
$ws.Range("M14").Value = 0.6068583333333333
$ws.Range("N14").Value = 1.820575
$ws.Range("O14").Value = 0.2510732750714712
$ws.Range("P14").Value = 0.2510732750714712
$ws.Range("Q14").Value = 0.1016267216472222
$ws.Range("R14").Value = 0.914640494825
$ws.Range("S14").Value = 0.0001978801627186375
$ws.Range("T14").Value = 0.0001978801627186375
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1674636666666667
$ws.Range("H15").Value = 0.502391
$ws.Range("I15").Value = 0.0007881370992683645
$ws.Range("J15").Value = 0.0007881370992683645
$ws.Range("O15").Value = 0.08210468103768082
$ws.Range("P15").Value = 0.08210468103768082
$ws.Range("Q15").Value = 0.03323344375611111
$ws.Range("R15").Value = 0.299100993805
$ws.Range("S15").Value = 0.00006470974514939205
$ws.Range("T15").Value = 0.00006470974514939205
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1674636666666667
$ws.Range("H16").Value = 0.502391
$ws.Range("I16").Value = 0.0007881370992683645
$ws.Range("J16").Value = 0.0007881370992683645
$ws.Range("M16").Value = 0.6741056666666667
$ws.Range("N16").Value = 2.022317
$ws.Range("O16").Value = 0.278895267936071
$ws.Range("P16").Value = 0.278895267936071
$ws.Range("Q16").Value = 0.1128882066607778
$ws.Range("R16").Value = 1.015993859947
$ws.Range("S16").Value = 0.0002198077074708083
$ws.Range("T16").Value = 0.0002198077074708083
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1674636666666667
$ws.Range("H17").Value = 0.502391
$ws.Range("I17").Value = 0.0007881370992683645
$ws.Range("J17").Value = 0.0007881370992683645
$ws.Range("M17").Value = 0.9376410000000001
$ws.Range("N17").Value = 2.812923
$ws.Range("O17").Value = 0.3879267759547769
$ws.Range("P17").Value = 0.387926775954777
$ws.Range("Q17").Value = 0.157020799877
$ws.Range("R17").Value = 1.413187198893
$ws.Range("S17").Value = 0.0003057394839295266
$ws.Range("T17").Value = 0.0003057394839295267
$ws.Range("G18").Value = 102.475225
$ws.Range("H18").Value = 307.425675
$ws.Range("I18").Value = 0.4822808922435293
$ws.Range("J18").Value = 0.4822808922435293
$ws.Range("M18").Value = 0.6068583333333333
$ws.Range("N18").Value = 1.820575
$ws.Range("O18").Value = 0.2510732750714712
$ws.Range("P18").Value = 0.2510732750714712
$ws.Range("Q18").Value = 62.18794425145832
$ws.Range("R18").Value = 559.6914982631249
$ws.Range("S18").Value = 0.1210878431199742
$ws.Range("T18").Value = 0.1210878431199742
$ws.Range("G19").Value = 102.475225
$ws.Range("H19").Value = 307.425675
$ws.Range("I19").Value = 0.4822808922435293
$ws.Range("J19").Value = 0.4822808922435293
$ws.Range("O19").Value = 0.08210468103768082
$ws.Range("P19").Value = 0.08210468103768082
$ws.Range("Q19").Value = 20.33637919329166
$ws.Range("R19").Value = 183.027412739625
$ws.Range("S19").Value = 0.03959751882822309
$ws.Range("T19").Value = 0.03959751882822309
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 102.475225
$ws.Range("H20").Value = 307.425675
$ws.Range("I20").Value = 0.4822808922435293
$ws.Range("J20").Value = 0.4822808922435293
$ws.Range("M20").Value = 0.6741056666666667
$ws.Range("N20").Value = 2.022317
$ws.Range("O20").Value = 0.278895267936071
$ws.Range("P20").Value = 0.278895267936071
$ws.Range("Q20").Value = 69.07912986544166
$ws.Range("R20").Value = 621.712168788975
$ws.Range("S20").Value = 0.1345058586627065
$ws.Range("T20").Value = 0.1345058586627065
$ws.Range("G21").Value = 102.475225
$ws.Range("H21").Value = 307.425675
$ws.Range("I21").Value = 0.4822808922435293
$ws.Range("J21").Value = 0.4822808922435293
$ws.Range("M21").Value = 0.9376410000000001
$ws.Range("N21").Value = 2.812923
$ws.Range("O21").Value = 0.3879267759547769
$ws.Range("P21").Value = 0.387926775954777
$ws.Range("Q21").Value = 96.08497244422499
$ws.Range("R21").Value = 864.7647519980248
$ws.Range("S21").Value = 0.1870896716326255
$ws.Range("T21").Value = 0.1870896716326255
